# Update "想去人数" (F column) figures for several rows in the
# "展览" and "全部类型" worksheets, reflecting newly scraped counts.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1117
$ws1.Range("F7").Value  = 61
$ws1.Range("F8").Value  = 11420
$ws1.Range("F9").Value  = 4325
$ws1.Range("F12").Value = 17
$ws1.Range("F13").Value = 2528
$ws1.Range("F18").Value = 502
$ws1.Range("F19").Value = 11283
$ws1.Range("F20").Value = 11167

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 1117
$ws4.Range("F7").Value  = 61
$ws4.Range("F8").Value  = 11420
$ws4.Range("F9").Value  = 4325
$ws4.Range("F12").Value = 17
$ws4.Range("F13").Value = 2528
$ws4.Range("F19").Value = 502
$ws4.Range("F20").Value = 11283
$ws4.Range("F21").Value = 11167
